$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume/Number and date range) ---
$a8 = $ws.Range("A8")
$a8r1 = $a8.Characters(21, 2)
$a8r1.Text = "28"
$a8r1b = $a8.Characters(21, 2)
$a8r1b.Font.Size = 10
$a8r1b.Font.Name = "Andale WT"

$c9 = $ws.Range("C9")
$c9r1 = $c9.Characters(27, 9)
$c9r1.Text = "7/7/2025"
$c9r1b = $c9.Characters(27, 8)
$c9r1b.Font.Size = 10
$c9r1b.Font.Name = "Andale WT"
$c9r2 = $c9.Characters(46, 8)
$c9r2.Text = "7/13/2025"
$c9r2b = $c9.Characters(46, 9)
$c9r2b.Font.Size = 10
$c9r2b.Font.Name = "Andale WT"

# --- Text-valued numeric-looking cells ("0" / "***.*") ---
# Source cell A14 already carries style s=13 (right-aligned text style)
# used as a format donor via PasteSpecial(formats) so the converted
# cells keep the same style id as the target diff.
$styleDonor = $ws.Range("A14")
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "0"
$styleDonor.Copy()
$c.PasteSpecial(-4122)
$c = $ws.Range("E15")
$c.NumberFormat = "@"
$c.Value = "***.*"
$styleDonor.Copy()
$c.PasteSpecial(-4122)
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "0"
$styleDonor.Copy()
$c.PasteSpecial(-4122)
$c = $ws.Range("E27")
$c.NumberFormat = "@"
$c.Value = "***.*"
$styleDonor.Copy()
$c.PasteSpecial(-4122)
$c = $ws.Range("C29")
$c.NumberFormat = "@"
$c.Value = "0"
$styleDonor.Copy()
$c.PasteSpecial(-4122)
$c = $ws.Range("C30")
$c.NumberFormat = "@"
$c.Value = "0"
$styleDonor.Copy()
$c.PasteSpecial(-4122)
$c = $ws.Range("C33")
$c.NumberFormat = "@"
$c.Value = "0"
$styleDonor.Copy()
$c.PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Plain numeric updates ---
$ws.Range("L14").Value = -83.333333333333
$ws.Range("G15").Value = 2
$ws.Range("H15").Value = 300
$ws.Range("I15").Value = 43
$ws.Range("K15").Value = 30.30303030303
$ws.Range("L15").Value = 95.454545454545
$ws.Range("M15").Value = 53.571428571428
$ws.Range("N15").Value = -10.416666666666
$ws.Range("C16").Value = 8
$ws.Range("D16").Value = 8
$ws.Range("E16").Value = 0
$ws.Range("F16").Value = 22
$ws.Range("G16").Value = 24
$ws.Range("H16").Value = -8.333333333333
$ws.Range("I16").Value = 160
$ws.Range("J16").Value = 169
$ws.Range("K16").Value = -5.325443786982
$ws.Range("L16").Value = -5.882352941176
$ws.Range("M16").Value = -23.444976076555
$ws.Range("N16").Value = -76.710334788937
$ws.Range("C17").Value = 21
$ws.Range("D17").Value = 20
$ws.Range("E17").Value = 5
$ws.Range("F17").Value = 100
$ws.Range("G17").Value = 69
$ws.Range("H17").Value = 44.927536231884
$ws.Range("I17").Value = 564
$ws.Range("J17").Value = 519
$ws.Range("K17").Value = 8.670520231213
$ws.Range("L17").Value = 9.514563106796
$ws.Range("M17").Value = 115.267175572519
$ws.Range("N17").Value = -9.032258064516
$ws.Range("C18").Value = 3
$ws.Range("D18").Value = 3
$ws.Range("E18").Value = 0
$ws.Range("F18").Value = 23
$ws.Range("G18").Value = 19
$ws.Range("H18").Value = 21.052631578947
$ws.Range("I18").Value = 181
$ws.Range("J18").Value = 159
$ws.Range("K18").Value = 13.836477987421
$ws.Range("L18").Value = -8.121827411167
$ws.Range("M18").Value = -40.849673202614
$ws.Range("N18").Value = -90.093048713738
$ws.Range("C19").Value = 35
$ws.Range("D19").Value = 25
$ws.Range("E19").Value = 40
$ws.Range("F19").Value = 119
$ws.Range("G19").Value = 119
$ws.Range("H19").Value = 0
$ws.Range("I19").Value = 723
$ws.Range("J19").Value = 864
$ws.Range("K19").Value = -16.319444444444
$ws.Range("L19").Value = -11.936662606577
$ws.Range("M19").Value = 44.6
$ws.Range("N19").Value = -15.63593932322
$ws.Range("C20").Value = 3
$ws.Range("D20").Value = 12
$ws.Range("E20").Value = -75
$ws.Range("F20").Value = 18
$ws.Range("G20").Value = 45
$ws.Range("H20").Value = -60
$ws.Range("I20").Value = 102
$ws.Range("J20").Value = 182
$ws.Range("K20").Value = -43.956043956044
$ws.Range("L20").Value = -54.464285714285
$ws.Range("M20").Value = -40.697674418604
$ws.Range("N20").Value = -96.001568012544
$ws.Range("C21").Value = 71
$ws.Range("D21").Value = 68
$ws.Range("E21").Value = 4.411764705882
$ws.Range("F21").Value = 290
$ws.Range("G21").Value = 279
$ws.Range("H21").Value = 3.942652329749
$ws.Range("I21").Value = 1775
$ws.Range("J21").Value = 1928
$ws.Range("K21").Value = -7.935684647302
$ws.Range("L21").Value = -9.484956654767
$ws.Range("M21").Value = 19.851451721809
$ws.Range("N21").Value = -73.122350090854
$ws.Range("D23").Value = 1
$ws.Range("E23").Value = 300
$ws.Range("F23").Value = 16
$ws.Range("G23").Value = 11
$ws.Range("H23").Value = 45.454545454545
$ws.Range("I23").Value = 84
$ws.Range("J23").Value = 60
$ws.Range("K23").Value = 40
$ws.Range("L23").Value = 0
$ws.Range("M23").Value = 147.058823529412
$ws.Range("C24").Value = 65
$ws.Range("D24").Value = 81
$ws.Range("E24").Value = -19.753086419753
$ws.Range("F24").Value = 276
$ws.Range("G24").Value = 295
$ws.Range("H24").Value = -6.440677966101
$ws.Range("I24").Value = 2230
$ws.Range("J24").Value = 2245
$ws.Range("K24").Value = -0.668151447661
$ws.Range("L24").Value = 0.495718792248
$ws.Range("M24").Value = 16.327595200834
$ws.Range("C25").Value = 42
$ws.Range("D25").Value = 34
$ws.Range("E25").Value = 23.529411764705
$ws.Range("F25").Value = 163
$ws.Range("G25").Value = 137
$ws.Range("H25").Value = 18.978102189781
$ws.Range("I25").Value = 1332
$ws.Range("J25").Value = 1208
$ws.Range("K25").Value = 10.264900662251
$ws.Range("L25").Value = 26.136363636363
$ws.Range("C26").Value = 36
$ws.Range("D26").Value = 32
$ws.Range("E26").Value = 12.5
$ws.Range("F26").Value = 167
$ws.Range("G26").Value = 176
$ws.Range("H26").Value = -5.113636363636
$ws.Range("I26").Value = 1032
$ws.Range("J26").Value = 1015
$ws.Range("K26").Value = 1.67487684729
$ws.Range("L26").Value = 8.975712777191
$ws.Range("M26").Value = -2.733270499528
$ws.Range("F27").Value = 10
$ws.Range("G27").Value = 4
$ws.Range("H27").Value = 150
$ws.Range("L27").Value = 45.945945945945
$ws.Range("C28").Value = 6
$ws.Range("D28").Value = 3
$ws.Range("E28").Value = 100
$ws.Range("F28").Value = 16
$ws.Range("G28").Value = 17
$ws.Range("H28").Value = -5.882352941176
$ws.Range("I28").Value = 113
$ws.Range("J28").Value = 104
$ws.Range("K28").Value = 8.653846153846
$ws.Range("L28").Value = -2.586206896551
$ws.Range("G29").Value = 1
$ws.Range("H29").Value = 100
$ws.Range("L29").Value = -75
$ws.Range("G30").Value = 1
$ws.Range("H30").Value = 100
$ws.Range("L30").Value = -72.222222222222
$ws.Range("G31").Value = 1
$ws.Range("H31").Value = 0
